# Assignment_13 edits — applies the diff by rebuilding the affected
# paragraphs' OOXML via Range.InsertXML (so we can emit <w:proofErr/>
# markers the Word OM itself has no property for).

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Get-ParaByText($doc, $substr) {
    $n = $doc.Paragraphs.Count
    for ($i = 1; $i -le $n; $i++) {
        $p = $doc.Paragraphs($i)
        if ($p.Range.Text -like "*$substr*") {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------
# 1) Answer 1: "Excel is superior ..." -> new Excel-vs-CSV answer.
# ---------------------------------------------------------------
$p = Get-ParaByText $d "Excel is superior to the CSV file format"
$xml = "<w:p $wns><w:pPr><w:spacing w:before='220'/><w:ind w:left='360'/></w:pPr>" +
       "<w:r><w:t xml:space='preserve'>Ans- </w:t></w:r>" +
       "<w:r><w:t>In Excel, spreadsheets can have values of datatypes other than strings; cells can have different fonts, sizes, or color settings; cells can have varying widths and heights; adjecent cells can be merged; and we can embed images and charts.</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 2) Question 2: split out "csv.reader" with proofErr gram markers.
# ---------------------------------------------------------------
$p = Get-ParaByText $d "2.What do you pass to"
$xml = "<w:p $wns><w:pPr><w:spacing w:before='220'/></w:pPr>" +
       "<w:r><w:t xml:space='preserve'>2.What do you pass to </w:t></w:r>" +
       "<w:proofErr w:type='gramStart'/>" +
       "<w:r><w:t>csv.reader</w:t></w:r>" +
       "<w:proofErr w:type='gramEnd'/>" +
       "<w:r><w:t>() and csv.writer() to create reader and writer objects?</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 3) Answer 2 block: 3 paragraphs -> 4 paragraphs.
#    ("Ans- First, call open..." / "This will create..." /
#     "to create a Writer object...")
# ---------------------------------------------------------------
$pStart = Get-ParaByText $d "First, call open"
$pEnd = Get-ParaByText $d "to create a Writer object"
$rng = $d.Range($pStart.Range.Start, $pEnd.Range.End)
$xml = (
    "<w:p $wns><w:pPr><w:spacing w:before='220'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Ans- </w:t></w:r>" +
    "<w:r><w:t xml:space='preserve'>We pass a file object, obtained from call to </w:t></w:r>" +
    "<w:proofErr w:type='gramStart'/>" +
    "<w:r><w:t>open(</w:t></w:r>" +
    "<w:proofErr w:type='gramEnd'/>" +
    "<w:r><w:t>).</w:t></w:r>" +
    "</w:p>" +
    "<w:p $wns><w:pPr><w:spacing w:before='220'/></w:pPr></w:p>" +
    "<w:p $wns><w:pPr><w:spacing w:before='220'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>This will create the object you can then pass to csv. </w:t></w:r>" +
    "<w:proofErr w:type='gramStart'/>" +
    "<w:r><w:t>writer(</w:t></w:r>" +
    "<w:proofErr w:type='gramEnd'/>" +
    "<w:r><w:t xml:space='preserve'>) </w:t></w:r>" +
    "</w:p>" +
    "<w:p $wns><w:pPr><w:spacing w:before='220'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>to create a Writer object. On Windows, you'll also need to pass a blank string for the </w:t></w:r>" +
    "<w:proofErr w:type='gramStart'/>" +
    "<w:r><w:t>open(</w:t></w:r>" +
    "<w:proofErr w:type='gramEnd'/>" +
    "<w:r><w:t>) function's newline keyword argument.</w:t></w:r>" +
    "</w:p>"
)
$rng.InsertXML($xml)

# ---------------------------------------------------------------
# 4) Answer 3: file-mode text swap (no proofErr here).
# ---------------------------------------------------------------
$p = Get-ParaByText $d "for reading"
$xml = "<w:p $wns><w:pPr><w:spacing w:before='220'/><w:ind w:left='360'/></w:pPr>" +
       "<w:r><w:t xml:space='preserve'>Ans- </w:t></w:r>" +
       "<w:r><w:t>File objects need to be opened in read-binary ('rb') for reader objects and write-binary ('wb') fro writer objects.</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 5) Answer 4: writerow() -> writerrow() w/ proofErr markers.
# ---------------------------------------------------------------
$p = Get-ParaByText $d "The most common method"
$xml = "<w:p $wns><w:pPr><w:spacing w:before='220'/><w:ind w:left='360'/></w:pPr>" +
       "<w:r><w:t>Ans-</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'>The </w:t></w:r>" +
       "<w:proofErr w:type='gramStart'/>" +
       "<w:r><w:t>writerrow(</w:t></w:r>" +
       "<w:proofErr w:type='gramEnd'/>" +
       "<w:r><w:t>) method.</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 6) Answer 6: loads() -> json.loads() w/ proofErr markers.
# ---------------------------------------------------------------
$p = Get-ParaByText $d "loads()"
$xml = "<w:p $wns><w:pPr><w:spacing w:before='220'/><w:ind w:left='360'/></w:pPr>" +
       "<w:r><w:t>Ans-</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:proofErr w:type='gramStart'/>" +
       "<w:r><w:t>json.loads</w:t></w:r>" +
       "<w:proofErr w:type='gramEnd'/>" +
       "<w:r><w:t>()</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

# ---------------------------------------------------------------
# 7) Answer 7: dumps() -> json.dumps() w/ proofErr markers.
# ---------------------------------------------------------------
$p = Get-ParaByText $d "dumps()"
$xml = "<w:p $wns><w:pPr><w:spacing w:before='220'/></w:pPr>" +
       "<w:r><w:t xml:space='preserve'>        Ans-</w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
       "<w:r><w:t xml:space='preserve'>  </w:t></w:r>" +
       "<w:proofErr w:type='gramStart'/>" +
       "<w:r><w:t>json.</w:t></w:r>" +
       "<w:r><w:t>dumps</w:t></w:r>" +
       "<w:proofErr w:type='gramEnd'/>" +
       "<w:r><w:t>()</w:t></w:r>" +
       "</w:p>"
$p.Range.InsertXML($xml)

Write-Output "done"
